$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto market data refresh (GitHub Actions scheduled update).
# Values that render as plain decimals (e.g. "1.00", "7.92") must be
# forced to Text format before assignment, otherwise Excel COM auto-
# converts them to numbers and mangles the original text formatting
# (trailing zeros, leading zeros, scientific notation, etc.).

$ws.Range("D2").Value = '60.599.88'
$ws.Range("E2").Value = '  -3.26%  '
$ws.Range("D3").Value = '3.346.05'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.99'
$ws.Range("E5").Value = '  -2.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.58'
$ws.Range("E6").Value = '  -1.11%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.484'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.92'
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("E10").Value = '  -1.40%  '
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("D12").Value = '3.915.20'
$ws.Range("E12").Value = '  -3.01%  '
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.72'
$ws.Range("E14").Value = '  -2.10%  '
$ws.Range("D15").Value = '3.339.75'
$ws.Range("E15").Value = '  -3.24%  '
$ws.Range("E16").Value = '  -1.71%  '
$ws.Range("D17").Value = '60.592.73'
$ws.Range("E17").Value = '  -3.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.28'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.55'
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.90'
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.88'
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.559'
$ws.Range("E22").Value = '  -0.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.79'
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '3.490.27'
$ws.Range("E25").Value = '  -2.55%  '
$ws.Range("E26").Value = '  -5.80%  '
$ws.Range("E27").Value = '  -4.57%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.34'
$ws.Range("E29").Value = '  -4.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.08'
$ws.Range("E30").Value = '  -1.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.69'
$ws.Range("E32").Value = '  -3.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.90'
$ws.Range("E33").Value = '  -1.44%  '
$ws.Range("E34").Value = '  -3.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.30'
$ws.Range("E35").Value = '  -1.18%  '
$ws.Range("E36").Value = '  -4.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.81'
$ws.Range("E37").Value = '  -2.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '167.09'
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.96'
$ws.Range("E39").Value = '  -12.32%  '
$ws.Range("D40").Value = '3.380.88'
$ws.Range("E40").Value = '  -2.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0747'
$ws.Range("E41").Value = '  -3.48%  '
$ws.Range("E42").Value = '  -3.65%  '
$ws.Range("E43").Value = '  -1.59%  '
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("E45").Value = '  -5.01%  '
$ws.Range("D46").Value = '2.454.36'
$ws.Range("E46").Value = '  -4.51%  '
$ws.Range("E47").Value = '  -3.37%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.36'
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("E50").Value = '  -2.00%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.816'
$ws.Range("E51").Value = '  +0.18%  '
